$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - copy the format from G1 (bold/border/centered header style)
# then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
